# "Fix in OnlyMuons column"
#
# Column K on Sheet1 is the "2p0E34_OnlyMuons" prescale column. It had been
# populated as a straight copy of the other menu columns, but this menu
# should carry a prescale of 0 (i.e. disabled) for all non-muon-only /
# non-special-case seeds. This restores K to 0 (highlighted in red, the
# workbook's existing convention for "forced-to-zero" prescales already
# seen on a handful of rows) for every affected seed row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rowsToZero = @(
    110,111,112,114,115,117,118,119,121,122,123,124,125,126,127,128,129,130,
    131,132,134,135,136,137,138,139,141,142,143,144,146,147,149,150,151,152,
    153,154,155,156,157,158,159,160,165,166,167,168,169,171,186,187,188,189,
    190,191,193,197,198,199,200,201,203,204,205,206,207,208,209,210,211,213,
    214,215,216,219,220,222,223,224,225,226,229,230,231,235,236,237,239,240,
    241,242,243,244,246,247,248,249,250,251,252,253,254,255,256,263,264,265,
    266,267,276,277,278,279,280,281,282,283,285,287,288,289,290,291,292,293,
    294,295,296,297,298,299,300,301,302,303,304,305,307,308,309,310,311,312,
    313,314,315,316,317,318,319,320,321,322,323,324,326,328,331,332,333,334,
    335,336,337,340,341,342,343,344
)

foreach ($r in $rowsToZero) {
    $cell = $ws.Range("K" + $r)
    $cell.Value = 0
    $cell.Font.Color = 255
}

# Reflect the resulting "select all" state the author saved with.
[void]$ws.Range("A1:M396").Select()

Write-Output ("Zeroed K column (OnlyMuons) for " + $rowsToZero.Count + " rows")
